$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (Primera/Segunda) for "Locoto" needs to be
# inserted before the existing row 110, pushing all subsequent rows (110-134)
# down by two rows (to 112-136). Insert two blank rows at 110:111 first.
$ws.Rows("110:111").Insert()

# Row 110 - new "Primera" observation
$ws.Range("A110").Value = 1
$ws.Range("B110").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C110").Value = "Arica y Parinacota"
$ws.Range("D110").Value = 44816
$ws.Range("E110").Value = 15
$ws.Range("F110").Value = 100112042
$ws.Range("G110").Value = "Locoto"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 120
$ws.Range("K110").Value = 17000
$ws.Range("L110").Value = 18000
$ws.Range("M110").Value = 17500
$ws.Range("N110").Value = "$/caja 20 kilos"
$ws.Range("O110").Value = "Región de Arica y Parinacota"
$ws.Range("P110").Value = 875
$ws.Range("Q110").Value = 20
$ws.Range("R110").Value = "Hortaliza"

# Row 111 - new "Segunda" observation
$ws.Range("A111").Value = 1
$ws.Range("B111").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C111").Value = "Arica y Parinacota"
$ws.Range("D111").Value = 44816
$ws.Range("E111").Value = 15
$ws.Range("F111").Value = 100112042
$ws.Range("G111").Value = "Locoto"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Segunda"
$ws.Range("J111").Value = 140
$ws.Range("K111").Value = 14000
$ws.Range("L111").Value = 15000
$ws.Range("M111").Value = 14500
$ws.Range("N111").Value = "$/caja 20 kilos"
$ws.Range("O111").Value = "Región de Arica y Parinacota"
$ws.Range("P111").Value = 725
$ws.Range("Q111").Value = 20
$ws.Range("R111").Value = "Hortaliza"
